$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "settings" - just move the active-cell selection
# ---------------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("settings")
$wsSettings.Activate()
$wsSettings.Range("B3").Select()

# ---------------------------------------------------------------------------
# Sheet "tier" - clear a stray boolean flag and drop the last (bad) row
# ---------------------------------------------------------------------------
$wsTier = $wb.Worksheets.Item("tier")
$wsTier.Activate()

# M2 was incorrectly flagged TRUE - clear it back to blank
$wsTier.Cells.Item(2, 13).Value = ""

# Row 7 (CHRMAUSER / CHRMAUSER_HJBB / t2-btc-2) was a bad acronym row - remove it
$wsTier.Rows("7:7").Delete()

$wsTier.Range("M14").Select()

# ---------------------------------------------------------------------------
# Sheet "ir" - move the TRUE flag from row 2 to row 6
# ---------------------------------------------------------------------------
$wsIr = $wb.Worksheets.Item("ir")
$wsIr.Activate()

$wsIr.Cells.Item(2, 5).Value = ""
$wsIr.Cells.Item(6, 5).Value = $true

$wsIr.Range("E3").Select()

# ---------------------------------------------------------------------------
# Sheet "simpleton" - replace the bad acronym rows with corrected data
# ---------------------------------------------------------------------------
$wsSimple = $wb.Worksheets.Item("simpleton")
$wsSimple.Activate()

# Rows 2 & 3 need the "last row of table" border formatting that row 5 already
# has, so copy its formats down before we rewrite the values.
$wsSimple.Range("A5:H5").Copy()
$wsSimple.Range("A2:H3").PasteSpecial(-4122)

$wsSimple.Cells.Item(3, 1).Value = "ARCONICTP"
$wsSimple.Cells.Item(3, 2).Value = "HJBT"
$wsSimple.Cells.Item(3, 3).Value = ""
$wsSimple.Cells.Item(3, 4).Value = ""
$wsSimple.Cells.Item(3, 5).Value = ""
$wsSimple.Cells.Item(3, 6).Value = "Flatbed"
$wsSimple.Cells.Item(3, 7).Value = ""
$wsSimple.Cells.Item(3, 8).Value = ""

$wsSimple.Cells.Item(2, 1).Value = "ARCONICTP"
$wsSimple.Cells.Item(2, 2).Value = "HJBB"
$wsSimple.Cells.Item(2, 3).Value = ""
$wsSimple.Cells.Item(2, 4).Value = ""
$wsSimple.Cells.Item(2, 5).Value = ""
$wsSimple.Cells.Item(2, 6).Value = "ICS"
$wsSimple.Cells.Item(2, 7).Value = ""
$wsSimple.Cells.Item(2, 8).Value = ""

# Old rows 4-6 (ALENDIST/BARCH4/... and TPCCOV/COCHDQ/... acronym junk) go away
$wsSimple.Rows("4:6").Delete()

$wsSimple.Range("H4").Select()

# ---------------------------------------------------------------------------
# Sheet "missingCode" - becomes the active tab, selection stays at A4
# ---------------------------------------------------------------------------
$wsMissing = $wb.Worksheets.Item("missingCode")
$wsMissing.Activate()
$wsMissing.Range("A4").Select()
